# Swap the order of "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
# in every cell of column G on the active sheet that contains that exact text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "dnasr281@gmail.com, System"
$newText = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G is the 7th column
    if ($cell.Value2 -eq $oldText) {
        $cell.Value2 = $newText
    }
}
